# Apply the "#5: fund, bonds, otherbonds, antique done" edit:
#  1. Remove the "其他有價證券" (other securities) sheet entirely - it held
#     garbage/mis-scraped header text rather than real data.
#  2. Fold the scraped "貝萊德世界礦業" fund-dealer row (which had spilled
#     into the wrong columns, with the numeric dealer id stuck in as text)
#     into the proper wide "name/owner/.../index" layout used by the other
#     property-type sheets, turning the old text "516257287" into the real
#     numeric total (5162572.87) and adding the trailing metadata columns.

$excel.DisplayAlerts = $false
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Delete the "其他有價證券" sheet (its rows were just mis-parsed
#    header/label text, not actual property records).
# ---------------------------------------------------------------------
$other = $wb.Worksheets.Item("其他有價證券")
$other.Delete()

# ---------------------------------------------------------------------
# 2) Rebuild "基金受益憑證" (fund) sheet with the full column layout.
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("基金受益憑證")

# -- header row 1: name / owner / dealer-name / quantity / face_value /
#    currency / total / property_category / category / date /
#    legislator_name / legislator_id / source_file / index
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "owner"
$ws.Range("D1").Value = "貝萊德世界礦業"
$ws.Range("E1").Value = "quantity"
$ws.Range("F1").Value = "face_value"
$ws.Range("G1").Value = "currency"
$ws.Range("H1").Value = "total"
$ws.Range("I1").Value = "property_category"
$ws.Range("J1").Value = "category"
$ws.Range("K1").Value = "date"
$ws.Range("L1").Value = "legislator_name"
$ws.Range("M1").Value = "legislator_id"
$ws.Range("N1").Value = "source_file"
$ws.Range("O1").Value = "index"

# -- data row 2: shift the old B/C/D values over, fix up D2/H2's types,
#    and append the record's property/category/legislator metadata.
$ws.Range("B2").Value = "永豐業銀行"
$ws.Range("C2").Value = "潘維剛"

# D2 must stay a literal text string "516257287" (not get silently
# re-interpreted as a number), so force text formatting before writing it.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "516257287"

# E2/F2/G2 (quantity/face_value/currency) are unchanged.
$ws.Range("H2").Value = 5162572.87
$ws.Range("I2").Value = "dealer"
$ws.Range("J2").Value = "normal"

# K2 (date) must stay literal text "2011-11-21", not become a date serial.
$ws.Range("K2").NumberFormat = "@"
$ws.Range("K2").Value = "2011-11-21"
$ws.Range("L2").Value = "潘維剛"
$ws.Range("M2").Value = 678
$ws.Range("N2").Value = "tmpcafb1"
$ws.Range("O2").Value = 78
